# Apply "Methylaction" -> "MethylAction" casing fix, update the date,
# and quote the lowercase "methylaction" package name reference.

$d = $word.ActiveDocument

# 1) Title: "Methylaction:" -> "MethylAction:"
$d.Content.Find.Execute("Methylaction:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MethylAction:", 2)

# 2) Date: "2015-05-07" -> "2015-05-25"
$d.Content.Find.Execute("2015-05-07", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2015-05-25", 2)

# 3) Remaining bare "Methylaction" occurrences -> "MethylAction"
#    (case-sensitive match so we don't disturb the lowercase "methylaction"
#    package-name references elsewhere in the document)
$d.Content.Find.Execute("Methylaction", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MethylAction", 2)

# 4) Quote the package name in "load the methylaction R package".
#    Use Find (no replacement) to locate the range, then assign Range.Text
#    directly so straight quotes are inserted verbatim (Find/Replace's
#    replacement text goes through AutoCorrect and turns " into smart
#    quotes, which we don't want here).
$rng = $d.Content
if ($rng.Find.Execute("load the methylaction R package into the session.")) {
    $rng.Text = 'load the "methylaction" R package into the session.'
}
